$d = $word.ActiveDocument

# 1. Change "ACTA No. 03" -> "ACTA No. 04" (the run with just "3")
$d.Content.Find.Execute("ACTA No. 03", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ACTA No. 04", 2) | Out-Null

# 2. Merge the multiple runs describing the documentation update into a single run
$d.Content.Find.Execute("Actualización y elaboración de Documentación (Creación del informe de retrospectiva, actualización documento especificación de requerimientos IEEE, actualización documento del Producto Backlog, creación del acta semanal)", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Actualización y elaboración de Documentación (Creación del informe de retrospectiva, actualización documento especificación de requerimientos IEEE, actualización documento del Producto Backlog, creación del acta semanal)", 2) | Out-Null
